# T5 work and minor stylesheet work
#
# Adds a new Hand (Hand 15) row to the "Hand " sheet and three new Glyph
# rows (g41, g42, g43) to the "Glyphs" sheet, along with their associated
# shared strings.

$wb = $excel.ActiveWorkbook

# --- "Hand " sheet: add row 15 -------------------------------------------
$wsHand = $wb.Worksheets.Item("Hand ")
$wsHand.Cells.Item(15, 1).Value = "Hand 15"
$wsHand.Cells.Item(15, 2).Value = "Unk. Hand in NLS Adv. MS. 72.1.33"
$wsHand.Cells.Item(15, 3).Value = "Transcription 5"

# --- "Glyphs" sheet: add rows 42-44 ---------------------------------------
$wsGlyphs = $wb.Worksheets.Item("Glyphs")

$wsGlyphs.Cells.Item(42, 1).Value = "g41"
$wsGlyphs.Cells.Item(42, 2).Value = "l with suspension stroke"

$wsGlyphs.Cells.Item(43, 1).Value = "g42"
$wsGlyphs.Cells.Item(43, 2).Value = "Superscript s"

$wsGlyphs.Cells.Item(44, 1).Value = "g43"
$wsGlyphs.Cells.Item(44, 2).Value = "Subscript i"

# --- Update the selection shown on the "Hand " sheet (C15) ----------------
# Select the new last cell so the saved view reflects it, then restore
# "Glyphs" as the active/visible sheet (it was tabSelected originally)
# and scroll/select it to match its new extent.
$wsHand.Range("C15").Select()

$wsGlyphs.Activate()
$excel.Goto($wsGlyphs.Range("A25"), $true)
$wsGlyphs.Range("B44").Select()
